# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) for the detail rows (16-60) is reversed:
# the newest-period-first ordering is flipped to an oldest-period-first
# ordering, mirroring the values top-to-bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 60

# Read current values in the Periodo Mora column (E)
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $values += , ($ws.Cells.Item($r, 5).Value())
}

# Write them back in reverse order (mirror top<->bottom)
$count = $values.Count
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 5).Value = $values[$count - 1 - $i]
}
